$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(0.282, 0.139, -0.013),
    @(0.076, 0.018, -0.075),
    @(0.029, -0.08400000000000001, -0.239),
    @(0.073, 0.322, -0.079),
    @(0.173, 0.039, -0.081),
    @(0.065, 0.016, -0.05),
    @(-0.004, -0.042, 0.599),
    @(-0.093, 0.298, 0.143),
    @(-0.107, -0.387, 0.029),
    @(-0.157, -0.185, -0.014)
)

for ($r = 0; $r -lt 10; $r++) {
    for ($c = 0; $c -lt 3; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}
